$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.358396053314209
$ws.Range("B1").Value = 2.927000045776367
$ws.Range("C1").Value = 6.026174068450928
$ws.Range("D1").Value = 2.138821601867676
$ws.Range("E1").Value = 0.766714870929718
